$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.418.70"
$ws.Range("E2").Value = "  +2.73%  "
$ws.Range("D3").Value = "3.461.55"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "578.16"
$ws.Range("E5").Value = "  +4.42%  "
$ws.Range("D6").Value = "158.38"
$ws.Range("E6").Value = "  +4.23%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.463.37"
$ws.Range("E8").Value = "  +3.50%  "
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").Value = "  +4.46%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  +6.27%  "
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").Value = "4.061.20"
$ws.Range("E13").Value = "  +3.80%  "
$ws.Range("D15").Value = "0.0000196"
$ws.Range("E15").Value = "  +9.07%  "
$ws.Range("D16").Value = "27.72"
$ws.Range("E16").Value = "  +3.45%  "
$ws.Range("D17").Value = "64.456.03"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "3.447.05"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "14.38"
$ws.Range("E20").Value = "  +5.06%  "
$ws.Range("D21").Value = "394.96"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("D24").Value = "72.91"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").Value = "  +28.53%  "
$ws.Range("D27").Value = "9.70"
$ws.Range("E27").Value = "  +9.85%  "
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "0.986"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "6.18"
$ws.Range("E30").Value = "  +11.53%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.04"
$ws.Range("E31").Value = "  +3.12%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.38"
$ws.Range("E32").Value = "  +7.47%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "6.64"
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").Value = "23.79"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("D36").Value = "7.06"
$ws.Range("E36").Value = "  +5.93%  "
$ws.Range("D37").Value = "160.90"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").Value = "0.0788"
$ws.Range("E39").Value = "  +7.09%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "27.54"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "2.934.68"
$ws.Range("D43").Value = "0.0322"
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "4.45"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.776"
$ws.Range("E45").Value = "  +3.23%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "42.18"
$ws.Range("E46").Value = "  +3.58%  "
$ws.Range("D47").Value = "24.04"
$ws.Range("E47").Value = "  +9.73%  "
$ws.Range("E48").Value = "  +5.14%  "
$ws.Range("D49").Value = "2.24"
$ws.Range("E49").Value = "  +26.50%  "
$ws.Range("D50").Value = "0.868"
$ws.Range("E50").Value = "  +7.93%  "
$ws.Range("D51").Value = "6.56"
$ws.Range("E51").Value = "  +4.64%  "
